$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.386.66'
$ws.Range('E2').Value = '  -0.44%  '
$ws.Range('D3').Value = '3.141.67'
$ws.Range('E3').Value = '  +1.13%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = "'533.15"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.05%  '
$ws.Range('D6').Value = "'143.02"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.54%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').Value = '3.140.23'
$ws.Range('E8').Value = '  +1.11%  '
$ws.Range('E9').Value = '  +0.77%  '
$ws.Range('E10').Value = '  -2.88%  '
$ws.Range('E11').Value = '  +0.34%  '
$ws.Range('E12').Value = '  +2.17%  '
$ws.Range('D13').Value = '3.682.25'
$ws.Range('E13').Value = '  +1.12%  '
$ws.Range('E14').Value = '  +3.28%  '
$ws.Range('E15').Value = '  -4.79%  '
$ws.Range('E16').Value = '  -0.38%  '
$ws.Range('D17').Value = '58.406.87'
$ws.Range('E17').Value = '  -0.44%  '
$ws.Range('D18').Value = '3.147.14'
$ws.Range('E18').Value = '  +1.29%  '
$ws.Range('D19').Value = "'6.12"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.37%  '
$ws.Range('D20').Value = "'12.87"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.56%  '
$ws.Range('D21').Value = "'8.00"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.54%  '
$ws.Range('D22').Value = "'344.29"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.75%  '
$ws.Range('E23').Value = '  -0.15%  '
$ws.Range('D24').Value = "'0.514"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.36%  '
$ws.Range('D25').Value = "'67.68"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.59%  '
$ws.Range('E26').Value = '  -0.74%  '
$ws.Range('D27').Value = "'1.00"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.11%  '
$ws.Range('D28').Value = '0.0₃0937'
$ws.Range('E28').Value = '  +1.88%  '
$ws.Range('D29').Value = "'7.50"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +3.27%  '
$ws.Range('E30').Value = '  -0.02%  '
$ws.Range('D31').Value = "'6.43"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -4.06%  '
$ws.Range('E32').Value = '  +1.40%  '
$ws.Range('D33').Value = "'21.14"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.59%  '
$ws.Range('D34').Value = "'1.20"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.99%  '
$ws.Range('B35').Value = 'Monero'
$ws.Range('C35').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D35').Value = "'158.33"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.83%  '
$ws.Range('B36').Value = 'NEARProtocol'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D36').Value = "'4.80"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.81%  '
$ws.Range('E37').Value = '  +2.36%  '
$ws.Range('D38').Value = "'26.37"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.29%  '
$ws.Range('D39').Value = "'1.26"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.52%  '
$ws.Range('E40').Value = '  +11.66%  '
$ws.Range('D41').Value = "'0.0671"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.16%  '
$ws.Range('E42').Value = '  +4.61%  '
$ws.Range('D43').Value = "'4.02"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.50%  '
$ws.Range('D44').Value = '3.182.20'
$ws.Range('E44').Value = '  +1.07%  '
$ws.Range('D45').Value = "'36.70"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.43%  '
$ws.Range('E46').Value = '  -0.03%  '
$ws.Range('E47').Value = '  +2.69%  '
$ws.Range('D48').Value = '2.282.15'
$ws.Range('E48').Value = '  -0.62%  '
$ws.Range('E49').Value = '  +4.02%  '
$ws.Range('D50').Value = "'20.72"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E51').Value = '  +1.56%  '
